$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the header row with two new values (14, 15) in P1 / Q1,
# matching the bold/centered/bordered style used by the rest of row 1.
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("P1:Q1").Font.Bold = $true
$ws.Range("P1:Q1").HorizontalAlignment = -4108
$ws.Range("P1:Q1").VerticalAlignment = -4160
$ws.Range("P1:Q1").Borders.LineStyle = 1

# Update existing data rows (2-25): swap the values in columns I, K, M, O
# and append new columns P and Q.
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 9).Value = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1  # K -> 1
    $ws.Cells.Item($r, 13).Value = 2  # M -> 2
    $ws.Cells.Item($r, 15).Value = 1  # O -> 1
    $ws.Cells.Item($r, 16).Value = 2  # P -> 2
    $ws.Cells.Item($r, 17).Value = 2  # Q -> 2
}
